$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TS Parameters")

$rows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,24,25,27,28,32,33,34,35,36,37,38,39,40,42,43,44,45,46,47,48,49,50,51,52,54,55,56,59,60)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).ClearContents()
}
